# table_config.xlsx adjustments: Vendor and Equipment table config changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert one new data row before the old row 13 so the 3 existing
#    "equipment" rows (old 13-15) and the trailing blank row (old 16)
#    shift down to rows 14-17, leaving row 13 free for a new
#    "S_ADDR_USAGE" vendor row.
# ------------------------------------------------------------------
$ws.Rows(13).Insert()

# ------------------------------------------------------------------
# 2. Normalize formatting for rows 13-17 (new row + shifted rows) so
#    they use the same look as the rest of the vendor/equipment rows
#    (copy format down from row 12).
# ------------------------------------------------------------------
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Re-write the full data block (rows 2-16) with the updated
#    content: ROWS column now uniformly 100, and a new
#    S_SUPPL_TAXNUMBERS row plus re-ordered / renamed rows.
# ------------------------------------------------------------------

# -- vendor --
$ws.Range("A2").Value = "vendor"
$ws.Range("B2").Value = "S_SUPPL_GEN"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "LIFNR, BU_GROUP, KTOKK, NAME_FIRST, NAME_FIRST_P, NAME_LAST_P, BPEXT, STREET, POST_CODE1, CITY1, COUNTRY, REGION, LANGU_CORR, TELNR_LONG, SMTP_ADDR"

$ws.Range("A3").Value = "vendor"
$ws.Range("B3").Value = "S_SUPPL_ADDR"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "LIFNR"

$ws.Range("A4").Value = "vendor"
$ws.Range("B4").Value = "S_SUPPL_COMPANY"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "LIFNR,BUKRS,AKONT,ZTERM1,ZWELS_01"

$ws.Range("A5").Value = "vendor"
$ws.Range("B5").Value = "S_SUPPL_PURCHASING"
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "LIFNR,EKORG,WAERS,ZTERM,INCO1,INCO2,KALKS,VSBED,WEBRE,KZAUT,BSTAE,KZRET"

$ws.Range("A6").Value = "vendor"
$ws.Range("B6").Value = "S_SUPP_BANK"
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "LIFNR,BANKS,BANKL,BANKN,IBAN,BKONT,BKREF,KOINH,EBPP_ACCNAME"

$ws.Range("A7").Value = "vendor"
$ws.Range("B7").Value = "S_SUPPL_PARTNER"
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = "LIFNR,EKORG,PARVW,LIFN2,DEFPA"

$ws.Range("A8").Value = "vendor"
$ws.Range("B8").Value = "S_SUPPL_TAXNUMBERS"
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = "TAXTYPE,TAXNUM"

$ws.Range("A9").Value = "vendor"
$ws.Range("B9").Value = "S_SUPPL_WITH_TAX"
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "LIFNR,BUKRS"

$ws.Range("A10").Value = "vendor"
$ws.Range("B10").Value = "S_LFA1_TEXT"
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "LIFNR,TDSPRAS"

$ws.Range("A11").Value = "vendor"
$ws.Range("B11").Value = "S_LFM1_TEXT"
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "LIFNR"

$ws.Range("A12").Value = "vendor"
$ws.Range("B12").Value = "S_ROLES"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "LIFNR,BP_ROLE"

$ws.Range("A13").Value = "vendor"
$ws.Range("B13").Value = "S_ADDR_USAGE"
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "LIFNR,ADR_KIND"

# -- equipment --
$ws.Range("A14").Value = "equipment"
$ws.Range("B14").Value = "S_EQUI"
$ws.Range("C14").Value = 100
$ws.Range("D14").Value = 12
$ws.Range("E14").Value = "EQUNR,NRANGE_IND,EQTYP,DATAB,EQKTX,BRGEW,GEWEI,MATNR,GERNR"

$ws.Range("A15").Value = "equipment"
$ws.Range("B15").Value = "S_IHPA"
$ws.Range("C15").Value = 100
$ws.Range("D15").Value = 13
$ws.Range("E15").Value = "EQUNR,PARVW"

$ws.Range("A16").Value = "equipment"
$ws.Range("B16").Value = "S_TEXTS_EQUI"
$ws.Range("C16").Value = 100
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = "EQUNR,SPRAS,TEXT_DESCR"

# row 17 stays blank (trailing blank row), matches old trailing row.

# ------------------------------------------------------------------
# 4. Row heights: a handful of rows differ from the default pattern.
# ------------------------------------------------------------------
$ws.Rows(3).RowHeight = 17.25
$ws.Rows(4).RowHeight = 17.25
$ws.Rows(13).RowHeight = 19.5
$ws.Rows(14).RowHeight = 17.25
$ws.Rows(15).RowHeight = 17.25
$ws.Rows(16).RowHeight = 17.25
$ws.Rows(17).RowHeight = 19.5

# ------------------------------------------------------------------
# 5. Column widths: column C is no longer the widest (it used to hold
#    the long combined column list), column B widened slightly.
# ------------------------------------------------------------------
$ws.Columns(2).ColumnWidth = 25.314523809523808
$ws.Columns(3).ColumnWidth = 19.314523809523808

$wb.Save()
